$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.4550319312332322
$ws.Range("C2").Value = 0.7785207984001278
$ws.Range("D2").Value = 0.7737883368274809
$ws.Range("E2").Value = 0.3901806420216862
$ws.Range("F2").Value = 0.5359095446178412
$ws.Range("G2").Value = 0.4929382061565302
$ws.Range("H2").Value = -168.6418104026258
$ws.Range("I2").Value = 347.2836208052515
$ws.Range("J2").Value = 364.6868154219615

# Row 3
$ws.Range("B3").Value = 0.4311398805418004
$ws.Range("C3").Value = 0.802576578552973
$ws.Range("D3").Value = 0.7974927136230067
$ws.Range("E3").Value = 0.3769268199153417
$ws.Range("F3").Value = 0.5700002847912971
$ws.Range("G3").Value = 0.5213210717488024
$ws.Range("H3").Value = -154.8579494692745
$ws.Range("I3").Value = 321.7158989385489
$ws.Range("J3").Value = 342.5997324786009

# Row 4
$ws.Range("B4").Value = 0.4069223654520485
$ws.Range("C4").Value = 0.8240148040955704
$ws.Range("D4").Value = 0.8187049059432816
$ws.Range("E4").Value = 0.3748515519166212
$ws.Range("F4").Value = 0.5776717687493665
$ws.Range("G4").Value = 0.5208198914656275
$ws.Range("H4").Value = -140.9921846603459
$ws.Range("I4").Value = 295.9843693206919
$ws.Range("J4").Value = 320.3488417840858

# Row 5
$ws.Range("B5").Value = 0.3863796922919697
$ws.Range("C5").Value = 0.8411693729025989
$ws.Range("D5").Value = 0.8356687451243339
$ws.Range("E5").Value = 0.3606162054705455
$ws.Range("F5").Value = 0.6131341330531762
$ws.Range("G5").Value = 0.5524492911791646
$ws.Range("H5").Value = -128.7841684182257
$ws.Range("I5").Value = 273.5683368364514
$ws.Range("J5").Value = 301.4134482231873

# Row 6
$ws.Range("B6").Value = 0.3658354831187558
$ws.Range("C6").Value = 0.8579625001802358
$ws.Range("D6").Value = 0.8524045110568537
$ws.Range("E6").Value = 0.3555229486032236
$ws.Range("F6").Value = 0.6302285823871159
$ws.Range("G6").Value = 0.5636697272167968
$ws.Range("H6").Value = -115.5993176653989
$ws.Range("I6").Value = 249.1986353307978
$ws.Range("J6").Value = 280.5243856408757

# Row 7
$ws.Range("B7").Value = 0.3461027877466166
$ws.Range("C7").Value = 0.8712176488721172
$ws.Range("D7").Value = 0.865593965416751
$ws.Range("E7").Value = 0.3489711867271373
$ws.Range("F7").Value = 0.6409659518932338
$ws.Range("G7").Value = 0.5676936971775672
$ws.Range("H7").Value = -103.8268112889457
$ws.Range("I7").Value = 227.6536225778914
$ws.Range("J7").Value = 262.4600118113113

# Row 8
$ws.Range("B8").Value = 0.3347057841674111
$ws.Range("C8").Value = 0.8803259477802978
$ws.Range("D8").Value = 0.8745521996468911
$ws.Range("E8").Value = 0.3358510955308506
$ws.Range("F8").Value = 0.6706055204292469
$ws.Range("G8").Value = 0.5951192855276159
$ws.Range("H8").Value = -95.00581615692917
$ws.Range("I8").Value = 212.0116323138583
$ws.Range("J8").Value = 250.2986604706202

# Row 9
$ws.Range("B9").Value = 0.3183557616374194
$ws.Range("C9").Value = 0.8915271850480235
$ws.Range("D9").Value = 0.88579293932369
$ws.Range("E9").Value = 0.3331767162507768
$ws.Range("F9").Value = 0.678530003201625
$ws.Range("G9").Value = 0.5964525572105506
$ws.Range("H9").Value = -83.18158892310285
$ws.Range("I9").Value = 190.3631778462057
$ws.Range("J9").Value = 232.1308449263096
